$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

$ws.Range("C2").Value  = "Mildly Glib"
$ws.Range("C3").Value  = "Mildly Grandiose"
$ws.Range("C4").Value  = "Mildly Conniving"
$ws.Range("C5").Value  = "Mildly Deceptive"
$ws.Range("C6").Value  = "Mildly Unremorseful"
$ws.Range("C7").Value  = "Mildly Callous"
$ws.Range("C8").Value  = "Mildly Inexpressive"
$ws.Range("C9").Value  = "Mildly Irresponsible"
$ws.Range("C10").Value = "Mildly Sensation Seeking"
$ws.Range("C11").Value = "Mildly Unrealistic"
$ws.Range("C12").Value = "Mildly Impulsive"
$ws.Range("C13").Value = "Mildly Irresponsible"
$ws.Range("C14").Value = "Mildly Parasitic"
$ws.Range("C15").Value = "Mildly Noncommittal"
$ws.Range("C16").Value = "Mildly Promiscuous"
$ws.Range("C17").Value = "Mildly Emotionally Controlled"
$ws.Range("C18").Value = "Mildly Problematic"
$ws.Range("C19").Value = "Mildly Delinquent"
$ws.Range("C20").Value = "Mildly Noncompliant"
